# Auto-generated edit script: updates leve-profit derived columns (H-N)
# per scheduled market-data refresh, across all 8 crafting-sheet tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 79.5
$ws.Range("I2").Value = 79.5
$ws.Range("K2").Value = 79.5
$ws.Range("M2").Value = 33.5
$ws.Range("H6").Value = 1263.5
$ws.Range("I6").Value = 407.55554
$ws.Range("K6").Value = 1222.66662
$ws.Range("M6").Value = -1110.66662
$ws.Range("H51").Value = 8793.5
$ws.Range("I51").Value = 8793.5
$ws.Range("K51").Value = 8793.5
$ws.Range("M51").Value = -8309.5
$ws.Range("H62").Value = 2999
$ws.Range("I62").Value = 2999
$ws.Range("K62").Value = 2999
$ws.Range("M62").Value = -2375
$ws.Range("H65").Value = 2999
$ws.Range("I65").Value = 2999
$ws.Range("K65").Value = 14995
$ws.Range("M65").Value = -11875
$ws.Range("H86").Value = 5633.7
$ws.Range("I86").Value = 3666.6667
$ws.Range("K86").Value = 3666.6667
$ws.Range("M86").Value = -2543.6667
$ws.Range("H89").Value = 5633.7
$ws.Range("I89").Value = 3666.6667
$ws.Range("K89").Value = 18333.3335
$ws.Range("M89").Value = -12717.3335
$ws.Range("H96").Value = 2515.6
$ws.Range("I96").Value = 1526.3334
$ws.Range("J96").Value = 3999.5
$ws.Range("K96").Value = 4579.0002
$ws.Range("L96").Value = 11998.5
$ws.Range("M96").Value = -3206.0002
$ws.Range("N96").Value = -14744.5
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H129").Value = 666.6667
$ws.Range("I129").Value = 500
$ws.Range("J129").Value = 1000
$ws.Range("K129").Value = 1500
$ws.Range("L129").Value = 3000
$ws.Range("M129").Value = 3500
$ws.Range("N129").Value = -13000
$ws.Range("H137").Value = 1024.1428
$ws.Range("I137").Value = 584.5
$ws.Range("K137").Value = 1753.5
$ws.Range("M137").Value = 796.5
$ws.Range("H138").Value = 2830.087
$ws.Range("J138").Value = 2964.7058
$ws.Range("L138").Value = 8894.117400000001
$ws.Range("N138").Value = -19174.1174

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1689.8667
$ws.Range("I61").Value = 1689.8667
$ws.Range("K61").Value = 1689.8667
$ws.Range("M61").Value = -1477.8667
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H88").Value = 2480.4
$ws.Range("I88").Value = 1824.75
$ws.Range("J88").Value = 2917.5
$ws.Range("K88").Value = 1824.75
$ws.Range("L88").Value = 2917.5
$ws.Range("M88").Value = -1418.75
$ws.Range("N88").Value = -3729.5
$ws.Range("H91").Value = 2480.4
$ws.Range("I91").Value = 1824.75
$ws.Range("J91").Value = 2917.5
$ws.Range("K91").Value = 1824.75
$ws.Range("L91").Value = 2917.5
$ws.Range("M91").Value = -420.75
$ws.Range("N91").Value = -5725.5
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H102").Value = 2362.25
$ws.Range("I102").Value = 2362.25
$ws.Range("K102").Value = 2362.25
$ws.Range("M102").Value = -740.25
$ws.Range("H110").Value = 571.1429000000001
$ws.Range("I110").Value = 499.66666
$ws.Range("K110").Value = 499.66666
$ws.Range("M110").Value = 1545.33334
$ws.Range("H136").Value = 1689.8667
$ws.Range("I136").Value = 1689.8667
$ws.Range("K136").Value = 5069.6001
$ws.Range("M136").Value = -2519.6001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2632.6667
$ws.Range("I86").Value = 3449
$ws.Range("K86").Value = 3449
$ws.Range("M86").Value = -2326
$ws.Range("H89").Value = 2632.6667
$ws.Range("I89").Value = 3449
$ws.Range("K89").Value = 17245
$ws.Range("M89").Value = -11629
$ws.Range("H99").Value = 3755.6
$ws.Range("I99").Value = 3755.6
$ws.Range("K99").Value = 3755.6
$ws.Range("M99").Value = -2257.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 650
$ws.Range("I13").Value = 650
$ws.Range("K13").Value = 650
$ws.Range("M13").Value = -511
$ws.Range("H31").Value = 2269.8
$ws.Range("I31").Value = 1599.6
$ws.Range("J31").Value = 2940
$ws.Range("K31").Value = 1599.6
$ws.Range("L31").Value = 2940
$ws.Range("M31").Value = -1304.6
$ws.Range("N31").Value = -3530
$ws.Range("H34").Value = 2269.8
$ws.Range("I34").Value = 1599.6
$ws.Range("J34").Value = 2940
$ws.Range("K34").Value = 1599.6
$ws.Range("L34").Value = 2940
$ws.Range("M34").Value = -1397.6
$ws.Range("N34").Value = -3344
$ws.Range("H62").Value = 1624.5
$ws.Range("I62").Value = 1499.5
$ws.Range("J62").Value = 1749.5
$ws.Range("K62").Value = 1499.5
$ws.Range("L62").Value = 1749.5
$ws.Range("M62").Value = -875.5
$ws.Range("N62").Value = -2997.5
$ws.Range("H65").Value = 1624.5
$ws.Range("I65").Value = 1499.5
$ws.Range("J65").Value = 1749.5
$ws.Range("K65").Value = 7497.5
$ws.Range("L65").Value = 8747.5
$ws.Range("M65").Value = -4377.5
$ws.Range("N65").Value = -14987.5
$ws.Range("H86").Value = 315237.12
$ws.Range("I86").Value = 188842.58
$ws.Range("K86").Value = 188842.58
$ws.Range("M86").Value = -187719.58
$ws.Range("H89").Value = 315237.12
$ws.Range("I89").Value = 188842.58
$ws.Range("K89").Value = 944212.8999999999
$ws.Range("M89").Value = -938596.8999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 368.625
$ws.Range("I8").Value = 368.625
$ws.Range("K8").Value = 1105.875
$ws.Range("M8").Value = -966.875
$ws.Range("H33").Value = 180
$ws.Range("I33").Value = 162.5
$ws.Range("K33").Value = 975
$ws.Range("M33").Value = -692
$ws.Range("H60").Value = 4997
$ws.Range("I60").Value = 4997
$ws.Range("K60").Value = 14991
$ws.Range("M60").Value = -14740
$ws.Range("H98").Value = 450
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H139").Value = 9995
$ws.Range("I139").Value = 9995
$ws.Range("K139").Value = 29985
$ws.Range("M139").Value = -24845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 3333.3333
$ws.Range("I59").Value = 2500
$ws.Range("K59").Value = 2500
$ws.Range("M59").Value = -1917
$ws.Range("H70").Value = 8983.666999999999
$ws.Range("I70").Value = 8983.666999999999
$ws.Range("K70").Value = 8983.666999999999
$ws.Range("M70").Value = -8713.666999999999
$ws.Range("H73").Value = 8983.666999999999
$ws.Range("I73").Value = 8983.666999999999
$ws.Range("K73").Value = 8983.666999999999
$ws.Range("M73").Value = -8047.666999999999
$ws.Range("H102").Value = 1586.1
$ws.Range("I102").Value = 1518.0555
$ws.Range("J102").Value = 2198.5
$ws.Range("K102").Value = 1518.0555
$ws.Range("L102").Value = 2198.5
$ws.Range("M102").Value = 103.9445000000001
$ws.Range("N102").Value = -5442.5
$ws.Range("H122").Value = 1755.3334
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 1266
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 3798
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -8698

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1200
$ws.Range("I7").Value = 1200
$ws.Range("K7").Value = 1200
$ws.Range("M7").Value = -1088
$ws.Range("H16").Value = 3255.7144
$ws.Range("I16").Value = 1465
$ws.Range("J16").Value = 14000
$ws.Range("K16").Value = 1465
$ws.Range("L16").Value = 14000
$ws.Range("M16").Value = -1295
$ws.Range("N16").Value = -14340
$ws.Range("H122").Value = 7532.4443
$ws.Range("I122").Value = 6827.5713
$ws.Range("J122").Value = 9999.5
$ws.Range("K122").Value = 20482.7139
$ws.Range("L122").Value = 29998.5
$ws.Range("M122").Value = -18032.7139
$ws.Range("N122").Value = -34898.5
$ws.Range("H126").Value = 1200
$ws.Range("I126").Value = 1200
$ws.Range("K126").Value = 3600
$ws.Range("M126").Value = -1130
$ws.Range("H136").Value = 4000
$ws.Range("I136").Value = 4000
$ws.Range("K136").Value = 12000
$ws.Range("M136").Value = -9450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 69950
$ws.Range("J46").Value = 69950
$ws.Range("L46").Value = 69950
$ws.Range("N46").Value = -70412
$ws.Range("H49").Value = 50000000
$ws.Range("J49").Value = 50000000
$ws.Range("L49").Value = 50000000
$ws.Range("N49").Value = -50000460
$ws.Range("H100").Value = 6972140
$ws.Range("I100").Value = 9958235
$ws.Range("K100").Value = 19916470
$ws.Range("M100").Value = -19915929
$ws.Range("H134").Value = 69950
$ws.Range("J134").Value = 69950
$ws.Range("L134").Value = 209850
$ws.Range("N134").Value = -214920
$ws.Range("H136").Value = 1030.0588
$ws.Range("I136").Value = 938.1875
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 2814.5625
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -264.5625
$ws.Range("N136").Value = -12600
